# Changes done to Update framework
#
# The "fulfiller" sample-data columns on row 4 (N4:AG4) of the test_data
# sheet are refreshed with a new, shorter set of placeholder names. The
# previous 19 distinct names (monika.negi, alka.bihal, ... deepti.bane)
# are replaced by a 4-value cycle (rahulhh, vijay, Nejha, perya) repeated
# across the N..AG range - this also naturally drops those 19 strings
# from the shared-strings table since they become unused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$values = @("rahulhh", "vijay", "Nejha", "perya")
$startCol = 14  # column N
$endCol   = 33  # column AG

for ($col = $startCol; $col -le $endCol; $col++) {
    $idx = ($col - $startCol) % 4
    $ws.Cells.Item(4, $col).Value = $values[$idx]
}

# Reflect the author's final on-screen scroll position / selection: the
# view had been scrolled so column W is the left-most visible column,
# with the active selection on AB17.
$window = $excel.ActiveWindow
$window.ScrollColumn = 23  # column W
$window.ScrollRow = 1
$ws.Range("AB17").Select()
